$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": revert Version / Date / Contact values ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.1.0"
$meta.Range("B8").Value = "2023-07-10T23:08:03+02:00"
$meta.Range("B10").Value = "No display for ContactDetail"

# --- Sheet "Include from FSIII": drop the 5 UUID concept rows that were added ---
$incl = $wb.Worksheets.Item("Include from FSIII")
$incl.Range("A2:A6").EntireRow.Delete()
